# Add a new game entry ("row 4") to the "Draymond Green" sheet, push the
# existing AVERAGE row down to row 5, and update the "final" sheet's
# references so they continue to point at the averages row.

$wb = $excel.ActiveWorkbook
$dray = $wb.Worksheets.Item("Draymond Green")
$final = $wb.Worksheets.Item("final")

# New game stats for Draymond Green (what used to be the AVERAGE row is now
# pushed down to row 5; row 4 becomes this new game's raw numbers).
$newGame = @(150, 5, 230, 43, 348, 65.22, 66.3, 69.9, 615)

for ($col = 1; $col -le 9; $col++) {
    $dray.Cells.Item(4, $col).Value = $newGame[$col - 1]
}

# Move the "promedios" label that lived in J4 down to J5.
$dray.Range("J4").Value = $null
$dray.Range("J5").Value = "promedios"

# Rebuild the AVERAGE formulas on row 5 (A5:I5) now that the source data is
# rows 2:4 instead of 2:3.
$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I")
foreach ($c in $cols) {
    $dray.Range($c + "5").Formula = "=AVERAGE(" + $c + "2:" + $c + "4)"
}

# Update "final" sheet row 4 formulas to pull from the new averages row
# (Draymond Green row 5 instead of row 4).
$finalCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
$drayCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt $finalCols.Length; $i++) {
    $final.Range($finalCols[$i] + "4").Formula = "='Draymond Green'!" + $drayCols[$i] + "5"
}

$excel.Calculate()
